# Add eos pairs (eos_eth, eos_usdt) to the binance exchange-limits sheet,
# following the existing pattern used by the other currency-pair rows
# (e.g. the etc_eth / etc_usdt rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: eos_eth -------------------------------------------------
$ws.Range("A13").Value = "eos_eth"
$ws.Range("B13").Value = 0.000001
$ws.Range("C13").Value = "ETH"
$ws.Range("D13").Value = 0.01
$ws.Range("E13").Value = "EOS"
$ws.Range("F13").Value = 0.01
$ws.Range("G13").Value = "EOS"
$ws.Range("H13").Value = 0.01
$ws.Range("I13").Value = "ETH"

# --- Row 14: eos_usdt --------------------------------------------------
$ws.Range("A14").Value = "eos_usdt"
$ws.Range("B14").Value = 0.0001
$ws.Range("C14").Value = "USDT"
$ws.Range("D14").Value = 0.01
$ws.Range("E14").Value = "EOS"
$ws.Range("F14").Value = 0.01
$ws.Range("G14").Value = "EOS"
$ws.Range("H14").Value = 10
$ws.Range("I14").Value = "USDT"

# Carry over the same cell formatting used by the rest of the table
# (the unit columns C/E/G/I carry the right-aligned/fill style) by
# copying formats down from the row above, the same way a user would
# fill these rows down from the existing table.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats

foreach ($col in @("C", "E", "G", "I")) {
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`13").PasteSpecial(-4122)
    $ws.Range("$col`14").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Column width tweaks observed on the sheet after the edit (column H is no
# longer auto "best fit" once a narrower value was set explicitly, and a
# few extra columns picked up explicit widths from the edit session).
$ws.Columns.Item(8).ColumnWidth = 14.5
$ws.Columns.Item(10).ColumnWidth = 9.65
$ws.Columns.Item(11).ColumnWidth = 14.65
$ws.Columns.Item(12).ColumnWidth = 14.65
$ws.Columns.Item(13).ColumnWidth = 11.2

# Move the active selection to the last edited cell, as Excel would after typing.
$ws.Range("I14").Select()
